$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$carrera = "Tecnicatura Universitaria en Procesamiento y Explotación de Datos"
$asignatura = "Espacio Integrador I"

# row 98
$ws.Cells.Item(98, 1).Value = $carrera
$ws.Cells.Item(98, 2).Value = $asignatura
$ws.Cells.Item(98, 3).Value = "La investigación cualitativa y el análisis computarizado de datos"
$ws.Cells.Item(98, 4).Value = "Echevarría, Hugo Darío"

# row 99 (autor entered before titulo)
$ws.Cells.Item(99, 1).Value = $carrera
$ws.Cells.Item(99, 2).Value = $asignatura
$ws.Cells.Item(99, 4).Value = "Scheaffer, Richard L"
$ws.Cells.Item(99, 3).Value = "elementos de muestreo"

# row 100
$ws.Cells.Item(100, 1).Value = $carrera
$ws.Cells.Item(100, 2).Value = $asignatura
$ws.Cells.Item(100, 3).Value = "tratamiento matemático de datos físico-químicos"
$ws.Cells.Item(100, 4).Value = "Spiridinov, V.P."

# row 101
$ws.Cells.Item(101, 1).Value = $carrera
$ws.Cells.Item(101, 2).Value = $asignatura
$ws.Cells.Item(101, 3).Value = "computadoras y procesamiento de datos"
$ws.Cells.Item(101, 4).Value = "Villanueva-Lara, Julio E"

# row 102
$ws.Cells.Item(102, 1).Value = $carrera
$ws.Cells.Item(102, 2).Value = $asignatura
$ws.Cells.Item(102, 3).Value = "Introducción a la programación y a las estructuras de datos"
$ws.Cells.Item(102, 4).Value = "Braunstein, Silvia L"

# row 103 (autor entered before titulo)
$ws.Cells.Item(103, 1).Value = $carrera
$ws.Cells.Item(103, 2).Value = $asignatura
$ws.Cells.Item(103, 4).Value = "Sheldom M. Ross"
$ws.Cells.Item(103, 3).Value = "introducción a la estadística"

# row 104
$ws.Cells.Item(104, 1).Value = $carrera
$ws.Cells.Item(104, 2).Value = $asignatura
$ws.Cells.Item(104, 3).Value = "Procesamiento de datos y análisis estadísticos"
$ws.Cells.Item(104, 4).Value = "Castañeda, Ma. Belén"

$ws.Range("A2:D104").Select()
$excel.ActiveWindow.ScrollRow = 91
